# Add data for 2023-12-02: updates several per-neighborhood worksheets
# plus the Citywide Totals / By Neighborhood rollup sheets with the
# incremented YTD violent-crime counts.
$wb = $excel.ActiveWorkbook

# Citywide Totals
$ws = $wb.Worksheets.Item(1)
$ws.Range("H2").Value = 105
$ws.Range("J2").Value = 120
$ws.Range("D3").Value = 133
$ws.Range("E3").Value = 143
$ws.Range("F3").Value = 135
$ws.Range("H3").Value = 153
$ws.Range("J3").Value = 227
$ws.Range("D4").Value = 8
$ws.Range("F4").Value = 8
$ws.Range("B6").Value = 371
$ws.Range("C6").Value = 473
$ws.Range("D6").Value = 411
$ws.Range("E6").Value = 465
$ws.Range("F6").Value = 523
$ws.Range("G6").Value = 433
$ws.Range("H6").Value = 437
$ws.Range("I6").Value = 497
$ws.Range("J6").Value = 410
$ws.Range("B7").Value = 496
$ws.Range("C7").Value = 627
$ws.Range("D7").Value = 642
$ws.Range("E7").Value = 689
$ws.Range("F7").Value = 758
$ws.Range("G7").Value = 662
$ws.Range("H7").Value = 711
$ws.Range("I7").Value = 828
$ws.Range("J7").Value = 778

# By Neighborhood
$ws = $wb.Worksheets.Item(2)
$ws.Range("F5").Value = 14
$ws.Range("B8").Value = 30
$ws.Range("E8").Value = 50
$ws.Range("F8").Value = 51
$ws.Range("J8").Value = 44
$ws.Range("H11").Value = 2
$ws.Range("B16").Value = 2
$ws.Range("D19").Value = 27
$ws.Range("F27").Value = 10
$ws.Range("B28").Value = 34
$ws.Range("D28").Value = 46
$ws.Range("H28").Value = 45
$ws.Range("F29").Value = 13
$ws.Range("E32").Value = 65
$ws.Range("D36").Value = 36
$ws.Range("F47").Value = 17
$ws.Range("F51").Value = 7
$ws.Range("C53").Value = 54
$ws.Range("D53").Value = 71
$ws.Range("F53").Value = 81
$ws.Range("H53").Value = 97
$ws.Range("I53").Value = 124
$ws.Range("J54").Value = 10
$ws.Range("D61").Value = 3
$ws.Range("E62").Value = 7
$ws.Range("D63").Value = 5
$ws.Range("E63").Value = 6
$ws.Range("C65").Value = 22
$ws.Range("F70").Value = 24
$ws.Range("H71").Value = 2
$ws.Range("I75").Value = 2
$ws.Range("F76").Value = 19
$ws.Range("G77").Value = 24
$ws.Range("J91").Value = 7
$ws.Range("H92").Value = 4
$ws.Range("G97").Value = 5
$ws.Range("B98").Value = 496
$ws.Range("C98").Value = 627
$ws.Range("D98").Value = 642
$ws.Range("E98").Value = 689
$ws.Range("F98").Value = 758
$ws.Range("G98").Value = 662
$ws.Range("H98").Value = 711
$ws.Range("I98").Value = 828
$ws.Range("J98").Value = 778

# Rogers Park
$ws = $wb.Worksheets.Item(3)
$ws.Range("F4").Value = 2
$ws.Range("F6").Value = 19

# Roseland
$ws = $wb.Worksheets.Item(4)
$ws.Range("G6").Value = 14
$ws.Range("G7").Value = 24

# Austin
$ws = $wb.Worksheets.Item(7)
$ws.Range("E3").Value = 7
$ws.Range("F3").Value = 6
$ws.Range("B5").Value = 21
$ws.Range("F5").Value = 36
$ws.Range("J5").Value = 24
$ws.Range("B6").Value = 30
$ws.Range("E6").Value = 50
$ws.Range("F6").Value = 51
$ws.Range("J6").Value = 44

# Garfield Park
$ws = $wb.Worksheets.Item(10)
$ws.Range("E6").Value = 52
$ws.Range("E7").Value = 65

# Grand Crossing
$ws = $wb.Worksheets.Item(11)
$ws.Range("D6").Value = 21
$ws.Range("D7").Value = 36

# Armour Square
$ws = $wb.Worksheets.Item(12)
$ws.Range("F3").Value = 1
$ws.Range("F6").Value = 14

# Englewood
$ws = $wb.Worksheets.Item(18)
$ws.Range("D4").Value = 1
$ws.Range("B6").Value = 30
$ws.Range("H6").Value = 26
$ws.Range("B7").Value = 34
$ws.Range("D7").Value = 46
$ws.Range("H7").Value = 45

# Loop
$ws = $wb.Worksheets.Item(22)
$ws.Range("H3").Value = 21
$ws.Range("C6").Value = 36
$ws.Range("D6").Value = 42
$ws.Range("F6").Value = 60
$ws.Range("I6").Value = 79
$ws.Range("C7").Value = 54
$ws.Range("D7").Value = 71
$ws.Range("F7").Value = 81
$ws.Range("H7").Value = 97
$ws.Range("I7").Value = 124

# West Loop
$ws = $wb.Worksheets.Item(24)
$ws.Range("J2").Value = 2
$ws.Range("J7").Value = 7

# West Pullman
$ws = $wb.Worksheets.Item(26)
$ws.Range("H2").Value = 1
$ws.Range("H5").Value = 4

# North Lawndale
$ws = $wb.Worksheets.Item(31)
$ws.Range("C5").Value = 18
$ws.Range("C6").Value = 22

# Near South Side
$ws = $wb.Worksheets.Item(35)
$ws.Range("E5").Value = 5
$ws.Range("E6").Value = 7

# Fuller Park
$ws = $wb.Worksheets.Item(36)
$ws.Range("F5").Value = 12
$ws.Range("F6").Value = 13

# New City
$ws = $wb.Worksheets.Item(43)
$ws.Range("D3").Value = 2
$ws.Range("E4").Value = 4
$ws.Range("D5").Value = 5
$ws.Range("E5").Value = 6

# Edgewater
$ws = $wb.Worksheets.Item(44)
$ws.Range("F4").Value = 6
$ws.Range("F5").Value = 10

# Wrigleyville
$ws = $wb.Worksheets.Item(45)
$ws.Range("G5").Value = 3
$ws.Range("G6").Value = 5

# Belmont Cragin
$ws = $wb.Worksheets.Item(46)
$ws.Range("G5").Value = 2
$ws.Range("G6").Value = 2

# Chatham
$ws = $wb.Worksheets.Item(47)
$ws.Range("D5").Value = 15
$ws.Range("D6").Value = 27

# Bucktown
$ws = $wb.Worksheets.Item(48)
$ws.Range("B4").Value = 2
$ws.Range("B5").Value = 2

# Little Village
$ws = $wb.Worksheets.Item(49)
$ws.Range("F4").Value = 6
$ws.Range("F5").Value = 7

# Portage Park
$ws = $wb.Worksheets.Item(58)
$ws.Range("E2").Value = 1
$ws.Range("E6").Value = 2

# Lower West Side
$ws = $wb.Worksheets.Item(60)
$ws.Range("J3").Value = 2
$ws.Range("J5").Value = 10

# Lake View
$ws = $wb.Worksheets.Item(65)
$ws.Range("F5").Value = 10
$ws.Range("F6").Value = 17

# Old Town
$ws = $wb.Worksheets.Item(67)
$ws.Range("F4").Value = 18
$ws.Range("F5").Value = 24

# Riverdale
$ws = $wb.Worksheets.Item(69)
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 2
